$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set cell A8 to "in progress" (new shared string)
$ws.Range("A8").Value = "in progress"

# Update the selection to D11
$ws.Range("D11").Select()
